# Actualización automática del mapa (2025-07-23 15:07:03)
# The "Recoleta / PEÑA 2079" case (row 63, Caso 6217) was removed from the
# dataset. Deleting the entire row shifts all subsequent rows up by one,
# which matches the rest of the diff (rows 64-69 data moving to 63-68 and
# the last row, 69, disappearing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(63).Delete() | Out-Null
